{"js": "// Update the answer text in each populated cell of the single table in the\n// document. The table has 20 rows x 5 columns; only every 4th row (0, 4, 8,\n// 12, 16) actually holds an answer like \"43\u00f72=21, 1\" \u2014 the rows in between\n// are blank spacer rows. We replace each populated cell's text with its new\n// value, matched by (row, col) position (not by old text, since some old\n// values repeat at different positions with different replacements).\n\nconst table = context.document.body.tables.getFirstOrNullObject();\ntable.load(\"rowCount\");\nawait context.sync();\n\nif (table.isNullObject) {\n  throw new Error(\"Expected a table in the document body, found none.\");\n}\n\nconst replacements = [\n  { row: 0, col: 0, before: \"43\u00f72=21, 1\", after: \"98\u00f78=12, 2\" },\n  { row: 0, col: 1, before: \"90\u00f72=45, 0\", after: \"64\u00f73=21, 1\" },\n  { row: 0, col: 2, before: \"96\u00f74=24, 0\", after: \"56\u00f75=11, 1\" },\n  { row: 0, col: 3, before: \"13\u00f75=2, 3\", after: \"86\u00f72=43, 0\" },\n  { row: 0, col: 4, before: \"84\u00f72=42, 0\", after: \"77\u00f76=12, 5\" },\n  { row: 4, col: 0, before: \"51\u00f74=12, 3\", after: \"87\u00f79=9, 6\" },\n  { row: 4, col: 1, before: \"37\u00f77=5, 2\", after: \"55\u00f76=9, 1\" },\n  { row: 4, col: 2, before: \"67\u00f74=16, 3\", after: \"53\u00f74=13, 1\" },\n  { row: 4, col: 3, before: \"91\u00f76=15, 1\", after: \"29\u00f74=7, 1\" },\n  { row: 4, col: 4, before: \"24\u00f79=2, 6\", after: \"74\u00f79=8, 2\" },\n  { row: 8, col: 0, before: \"17\u00f74=4, 1\", after: \"46\u00f74=11, 2\" },\n  { row: 8, col: 1, before: \"46\u00f77=6, 4\", after: \"46\u00f76=7, 4\" },\n  { row: 8, col: 2, before: \"56\u00f73=18, 2\", after: \"24\u00f78=3, 0\" },\n  { row: 8, col: 3, before: \"31\u00f72=15, 1\", after: \"14\u00f77=2, 0\" },\n  { row: 8, col: 4, before: \"61\u00f76=10, 1\", after: \"23\u00f78=2, 7\" },\n  { row: 12, col: 0, before: \"26\u00f72=13, 0\", after: \"81\u00f76=13, 3\" },\n  { row: 12, col: 1, before: \"68\u00f74=17, 0\", after: \"13\u00f74=3, 1\" },\n  { row: 12, col: 2, before: \"65\u00f76=10, 5\", after: \"42\u00f79=4, 6\" },\n  { row: 12, col: 3, before: \"23\u00f72=11, 1\", after: \"90\u00f76=15, 0\" },\n  { row: 12, col: 4, before: \"86\u00f73=28, 2\", after: \"45\u00f78=5, 5\" },\n  { row: 16, col: 0, before: \"46\u00f72=23, 0\", after: \"49\u00f79=5, 4\" },\n  { row: 16, col: 1, before: \"28\u00f74=7, 0\", after: \"50\u00f79=5, 5\" },\n  { row: 16, col: 2, before: \"85\u00f72=42, 1\", after: \"75\u00f73=25, 0\" },\n  { row: 16, col: 3, before: \"31\u00f72=15, 1\", after: \"86\u00f75=17, 1\" },\n  { row: 16, col: 4, before: \"16\u00f72=8, 0\", after: \"42\u00f75=8, 2\" },\n];\n\n// Load current values of every target cell first, so we can sanity-check\n// before mutating anything.\nconst cells = replacements.map((r) => table.getCell(r.row, r.col));\ncells.forEach((c) => c.load(\"value\"));\nawait context.sync();\n\nfor (let i = 0; i < replacements.length; i++) {\n  const expected = replacements[i].before;\n  const actual = cells[i].value;\n  if (actual !== expected) {\n    throw new Error(\n      `Cell (${replacements[i].row}, ${replacements[i].col}) text mismatch: ` +\n        `expected \"${expected}\" but found \"${actual}\"`\n    );\n  }\n}\n\nfor (let i = 0; i < replacements.length; i++) {\n  cells[i].value = replacements[i].after;\n}\n\nawait context.sync();\n", "ps1": "# Update the answer text in each populated cell of the single table in the\n# document. The table has 20 rows x 5 columns; only every 4th row (1, 5, 9,\n# 13, 17 in 1-based COM indexing) actually holds an answer like\n# \"43\u00f72=21, 1\" \u2014 the rows in between are blank spacer rows. We replace each\n# populated cell's text with its new value, matched by (row, col) position\n# (not by old text, since some old values repeat at different positions\n# with different replacements).\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$replacements = @(\n    @{ Row = 1;  Col = 1; Before = \"43\u00f72=21, 1\";  After = \"98\u00f78=12, 2\" },\n    @{ Row = 1;  Col = 2; Before = \"90\u00f72=45, 0\";  After = \"64\u00f73=21, 1\" },\n    @{ Row = 1;  Col = 3; Before = \"96\u00f74=24, 0\";  After = \"56\u00f75=11, 1\" },\n    @{ Row = 1;  Col = 4; Before = \"13\u00f75=2, 3\";   After = \"86\u00f72=43, 0\" },\n    @{ Row = 1;  Col = 5; Before = \"84\u00f72=42, 0\";  After = \"77\u00f76=12, 5\" },\n    @{ Row = 5;  Col = 1; Before = \"51\u00f74=12, 3\";  After = \"87\u00f79=9, 6\" },\n    @{ Row = 5;  Col = 2; Before = \"37\u00f77=5, 2\";   After = \"55\u00f76=9, 1\" },\n    @{ Row = 5;  Col = 3; Before = \"67\u00f74=16, 3\";  After = \"53\u00f74=13, 1\" },\n    @{ Row = 5;  Col = 4; Before = \"91\u00f76=15, 1\";  After = \"29\u00f74=7, 1\" },\n    @{ Row = 5;  Col = 5; Before = \"24\u00f79=2, 6\";   After = \"74\u00f79=8, 2\" },\n    @{ Row = 9;  Col = 1; Before = \"17\u00f74=4, 1\";   After = \"46\u00f74=11, 2\" },\n    @{ Row = 9;  Col = 2; Before = \"46\u00f77=6, 4\";   After = \"46\u00f76=7, 4\" },\n    @{ Row = 9;  Col = 3; Before = \"56\u00f73=18, 2\";  After = \"24\u00f78=3, 0\" },\n    @{ Row = 9;  Col = 4; Before = \"31\u00f72=15, 1\";  After = \"14\u00f77=2, 0\" },\n    @{ Row = 9;  Col = 5; Before = \"61\u00f76=10, 1\";  After = \"23\u00f78=2, 7\" },\n    @{ Row = 13; Col = 1; Before = \"26\u00f72=13, 0\";  After = \"81\u00f76=13, 3\" },\n    @{ Row = 13; Col = 2; Before = \"68\u00f74=17, 0\";  After = \"13\u00f74=3, 1\" },\n    @{ Row = 13; Col = 3; Before = \"65\u00f76=10, 5\";  After = \"42\u00f79=4, 6\" },\n    @{ Row = 13; Col = 4; Before = \"23\u00f72=11, 1\";  After = \"90\u00f76=15, 0\" },\n    @{ Row = 13; Col = 5; Before = \"86\u00f73=28, 2\";  After = \"45\u00f78=5, 5\" },\n    @{ Row = 17; Col = 1; Before = \"46\u00f72=23, 0\";  After = \"49\u00f79=5, 4\" },\n    @{ Row = 17; Col = 2; Before = \"28\u00f74=7, 0\";   After = \"50\u00f79=5, 5\" },\n    @{ Row = 17; Col = 3; Before = \"85\u00f72=42, 1\";  After = \"75\u00f73=25, 0\" },\n    @{ Row = 17; Col = 4; Before = \"31\u00f72=15, 1\";  After = \"86\u00f75=17, 1\" },\n    @{ Row = 17; Col = 5; Before = \"16\u00f72=8, 0\";   After = \"42\u00f75=8, 2\" }\n)\n\nforeach ($rep in $replacements) {\n    $cell = $t.Cell($rep.Row, $rep.Col)\n    $cellRange = $cell.Range\n    # Range.Text for a table cell includes the trailing cell-mark character;\n    # strip any trailing control characters before comparing.\n    $current = $cellRange.Text.TrimEnd([char]7, [char]13)\n    if ($current -ne $rep.Before) {\n        throw \"Cell ($($rep.Row), $($rep.Col)) text mismatch: expected '$($rep.Before)' but found '$current'\"\n    }\n    $cellRange.Text = $rep.After\n}\n"}
